# Nalco aluminium-ingot price table: a new circular was published, so a new
# "latest price" row is inserted right below the header and every existing
# row shifts down by one (Sl.no. values are bumped accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing data rows down by one to make room for the new entry.
$ws.Rows("2:2").Insert()

# E2/F2 hold a dd-mm-yyyy string and a URL respectively; force them to stay
# plain text so Excel doesn't auto-convert the date-like string into a date
# serial number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "09-10-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf"

$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 282.85

# Match the formatting of the rest of the table (Insert() only blanked the
# row, it copied the header's bold style by default).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)  # xlPasteFormats

# Row-insert shifts cell values/styles automatically but leaves the
# worksheet's <hyperlinks> ref-to-target mapping untouched, which would
# leave every link pointing one row off. Rebuild all of them from scratch
# against the now-correct F-column text.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $target = $cell.Value2
    $ws.Hyperlinks.Add($cell, $target) | Out-Null
}

# Hyperlinks.Add applies Excel's default hyperlink look (blue/underline);
# restore the table's normal centered style used by every other column.
$ws.Range("C3").Copy()
$ws.Range("F2:F12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A1").Select()
